$wb = $excel.ActiveWorkbook

# Update the Date value on the Metadata sheet
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2022-07-26T14:00:32+00:00"

# Add a new concept row on the Concepts sheet, copying formatting from the row above
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("A26:D26").Copy($concepts.Range("A27:D27"))
$concepts.Range("B27").Value = "PHEN"
$concepts.Range("C27").Value = "Clinical Sign (HPO)"
